# Add 20 new data rows (44-63) to the sheet, matching the COM automation
# behaviour: "impresion calidad y ajuste de . , en campo ancho"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @('120', '120', '120', '120', '122', '1', '221', '25/08/2025 10:40', 'A', '02', 'COVERING'),
    @('120', '120', '120', '120', '122', '2', '221', '25/08/2025 10:41', 'A', '02', 'COVERING'),
    @('120', '120', '120', '120', '122', '3', '221', '25/08/2025 10:43', 'A', '02', 'COVERING'),
    @('120', '120', '120', '120', '122', '4', '221', '25/08/2025 10:48', 'A', '02', 'COVERING'),
    @('120', '120', '120', '120', '122', '5', '221', '25/08/2025 10:48', 'A', '02', 'COVERING'),
    @('60', '121', '120', '120', '12', '1', '1211', '25/08/2025 10:51', 'A', '03', 'L.BLANCO'),
    @('120', '120', '150', '511', '15444', '1', '5454', '25/08/2025 10:55', 'A', '04', 'CART.GRIS'),
    @('120', '120', '150', '545', '4544', '1', '4444', '25/08/2025 11:00', 'A', '03', 'L.BLANCO'),
    @('120', '120', '150', '545', '4544', '2', '4444', '25/08/2025 11:03', 'A', '03', 'L.BLANCO'),
    @('120', '120', '150', '560', '54545', '1', '21215', '25/08/2025 11:29', 'A', '03', 'L.BLANCO'),
    @('120', '120', '150', '544', '1245', '1', '45455', '26/08/2025 08:33', 'A', '03', 'L.BLANCO'),
    @('130', '120', '150', '544', '12154', '1', '4545', '26/08/2025 08:36', 'B', '03', 'L.BLANCO'),
    @('120', '120', '150', '254', '5455', '1', '4455', '26/08/2025 08:38', 'B', '05', 'CART.BLANCA'),
    @('120', '120', '150', '544', '45454', '1', '1111', '26/08/2025 08:40', 'B', '05', 'CART.BLANCA'),
    @('120', '150', '150', '1501', '2222', '1', '22454', '26/08/2025 08:42', 'A', '01', 'ONDA LINER'),
    @('120', '120', '150', '545', '45578', '1', '88877', '26/08/2025 08:44', 'B', '05', 'CART.BLANCA'),
    @('120', '120', '150', '545', '45578', '2', '88877', '26/08/2025 08:44', 'B', '01', 'ONDA LINER'),
    @('62,5', '120', '150', '2455', '4455', '1', '5544', '26/08/2025 08:47', 'C', '05', 'CART.BLANCA'),
    @('62,5', '120', '170', '5444', '5454', '1', '54545', '26/08/2025 08:48', 'B', '05', 'CART.BLANCA'),
    @('62,5', '120', '170', '5444', '5454', '2', '54545', '26/08/2025 08:49', 'B', '01', 'ONDA LINER')
)

$startRow = 44
$endRow = $startRow + $newRows.Count - 1

# Force the destination range to be stored as text so values such as
# "02" (leading zero) and "62,5" (comma) are preserved exactly as typed,
# matching the inline-string cells already used throughout the sheet.
$rangeAddress = "A" + $startRow + ":K" + $endRow
$destRange = $ws.Range($rangeAddress)
$destRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}

Write-Host "Inserted rows $startRow to $endRow"
